$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Paragraph 1 ("**ID__AFFARS_5325_topic_10__ID**") gets the same paragraph
# border / indent treatment already used by the body paragraphs further
# down in the document: a border on all four edges with no line (just a
# w:space of 5), and the left indent bumped from 120 to 225 twips.
# ---------------------------------------------------------------------------
$p = $d.Paragraphs(1)
$pf = $p.Range.ParagraphFormat
$pf.Borders.DistanceFromTop = 5
$pf.Borders.DistanceFromLeft = 5
$pf.Borders.DistanceFromBottom = 5
$pf.Borders.DistanceFromRight = 5
$pf.LeftIndent = 225 * 0.05   ; # LeftIndent is in points; 225 twips = 11.25pt

# ---------------------------------------------------------------------------
# The paragraph currently holds two runs: the id placeholder text, then a
# run containing a single trailing space. Locate the id text, drop the
# trailing-space run entirely, and repoint the id placeholder at the new
# topic id.
# ---------------------------------------------------------------------------
$idRange = $d.Content
$found = $idRange.Find.Execute("**ID__AFFARS_5325_topic_10__ID**", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$p = $d.Paragraphs(1)
$paraMarkStart = $p.Range.End - 1   ; # position just before the paragraph mark

$trailing = $d.Range($idRange.End, $paraMarkStart)
if ($trailing.Start -lt $trailing.End) {
    $trailing.Delete()
}

$idRange.Text = "**ID__AFFARS_5325_603__ID**"
